$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "58×94=5452"
$t.Cell(1, 2).Range.Text = "21×29=609"
$t.Cell(1, 3).Range.Text = "12×11=132"
$t.Cell(1, 4).Range.Text = "40×72=2880"
$t.Cell(1, 5).Range.Text = "65×36=2340"

$t.Cell(5, 1).Range.Text = "11×15=165"
$t.Cell(5, 2).Range.Text = "18×83=1494"
$t.Cell(5, 3).Range.Text = "98×83=8134"
$t.Cell(5, 4).Range.Text = "72×67=4824"
$t.Cell(5, 5).Range.Text = "39×26=1014"

$t.Cell(10, 1).Range.Text = "95×76=7220"
$t.Cell(10, 2).Range.Text = "34×49=1666"
$t.Cell(10, 3).Range.Text = "45×64=2880"
$t.Cell(10, 4).Range.Text = "35×19=665"
$t.Cell(10, 5).Range.Text = "13×20=260"

$t.Cell(15, 1).Range.Text = "22×67=1474"
$t.Cell(15, 2).Range.Text = "44×39=1716"
$t.Cell(15, 3).Range.Text = "88×49=4312"
$t.Cell(15, 4).Range.Text = "12×72=864"
$t.Cell(15, 5).Range.Text = "59×76=4484"

$t.Cell(20, 1).Range.Text = "36×46=1656"
$t.Cell(20, 2).Range.Text = "53×38=2014"
$t.Cell(20, 3).Range.Text = "32×52=1664"
$t.Cell(20, 4).Range.Text = "23×92=2116"
$t.Cell(20, 5).Range.Text = "36×91=3276"
